# Apply the "parallel structure" data-processing rename/metadata edit:
#   - "1880s"  -> "1880Survey"
#   - "1940s"  -> "1940Survey"
#   - add a new "1880Metadata" sheet (ToDo-list style single note) at the end
#   - update sheet selections to match the author's saved state

$wb = $excel.ActiveWorkbook

# --- Rename the two data tabs ---------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "1880Survey"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "1940Survey"

# Put the selection on 1940Survey first so activating 1880Survey afterwards
# leaves it as the final/active (tabSelected) sheet, matching the source file.
$ws3.Activate()
[void]$ws3.Range("A4").Select()

# --- Add the new metadata sheet at the end of the workbook -----------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "1880Metadata"
$newSheet.Range("A1").Value = "Actually surveyed in 1881; tab reads 1880 for consistency in the data processing R script"

# --- Restore 1880Survey as the active/selected sheet ------------------------
$ws1.Activate()
[void]$ws1.Range("A2").Select()
